# RippleTemplate_Combination.xlsx update
#  - add a new "Assay" worksheet (settings/value table) at the end of the workbook
#  - make "Patterns" the active/selected sheet again (it had been "Compounds")
#  - update selections on "Patterns" and the new "Assay" sheet

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Add the new "Assay" worksheet after the last existing sheet (Barcodes)
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$assay = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$assay.Name = "Assay"

$assay.Range("A1").Value = "Setting"
$assay.Range("B1").Value = "Value"

$assay.Range("A2").Value = "DMSO Tolerance"
$assay.Range("B2").Value = 0.005

$assay.Range("A3").Value = "Well Volume (µL)"
$assay.Range("B3").Value = 25

$assay.Range("A4").Value = "Backfill (µL)"
$assay.Range("B4").Value = 10

$assay.Range("A5").Value = "Allowed Error"
$assay.Range("B5").Value = 0.1

$assay.Range("A6").Value = "Destination Replicates"
$assay.Range("B6").Value = 1

$assay.Range("A7").Value = "Use Intermediate Plates"
$assay.Range("B7").Value = 1

$assay.Range("A8").Value = "DMSO Normalization"
$assay.Range("B8").Value = 1

# Leave the new sheet's own selection spanning its whole settings table
# (this also makes it the tab that's current right after creation).
$null = $assay.Range("A1:B8").Select()

# ---------------------------------------------------------------------------
# 2) Switch back to "Patterns" as the active sheet/tab and set its selection
#    (was previously "Compounds"; doing this last is what makes "Patterns"
#    the tab that is persisted as active/selected on save).
# ---------------------------------------------------------------------------
$patterns = $wb.Worksheets.Item("Patterns")
$patterns.Activate()
$null = $patterns.Range("O10").Select()
